$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.226.17'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '2.191.75'
$ws.Range("E3").Value = '  -5.48%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '295.27'
$ws.Range("E5").Value = '  -3.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.32'
$ws.Range("E6").Value = '  -2.29%  '
$ws.Range("E7").Value = '  -2.79%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.465'
$ws.Range("E9").Value = '  -3.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0768'
$ws.Range("E10").Value = '  -4.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '29.02'
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.09'
$ws.Range("E12").Value = '  -10.04%  '
$ws.Range("E13").Value = '  -1.93%  '
$ws.Range("D14").Value = '2.527.76'
$ws.Range("E14").Value = '  -5.63%  '
$ws.Range("E15").Value = '  -1.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.91'
$ws.Range("E16").Value = '  -4.96%  '
$ws.Range("D17").Value = '2.187.25'
$ws.Range("E17").Value = '  -6.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.707'
$ws.Range("E18").Value = '  -4.00%  '
$ws.Range("D19").Value = '39.104.71'
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("E20").Value = '  -2.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.69'
$ws.Range("E21").Value = '  -5.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.77'
$ws.Range("E22").Value = '  -3.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.21'
$ws.Range("E23").Value = '  -2.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '224.91'
$ws.Range("E24").Value = '  -3.19%  '
$ws.Range("E26").Value = '  -4.83%  '
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("E28").Value = '  -2.59%  '
$ws.Range("E29").Value = '  -1.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.03'
$ws.Range("E30").Value = '  -0.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '149.47'
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.57'
$ws.Range("E32").Value = '  -6.30%  '
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.79'
$ws.Range("E34").Value = '  -5.23%  '
$ws.Range("E35").Value = '  -3.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0690'
$ws.Range("E36").Value = '  -3.31%  '
$ws.Range("E37").Value = '  -2.18%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.29'
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0960'
$ws.Range("E39").Value = '  -2.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.61'
$ws.Range("E40").Value = '  -3.95%  '
$ws.Range("E41").Value = '  -2.65%  '
$ws.Range("E42").Value = '  -4.26%  '
$ws.Range("D43").Value = '1.889.24'
$ws.Range("E43").Value = '  -2.18%  '
$ws.Range("E44").Value = '  -9.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0259'
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.04'
$ws.Range("E46").Value = '  -7.38%  '
$ws.Range("E47").Value = '  -2.96%  '
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '71.03'
$ws.Range("E49").Value = '  +1.36%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.398.26'
$ws.Range("E50").Value = '  -6.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.64'
$ws.Range("E51").Value = '  -4.89%  '
